$d = $word.ActiveDocument

function Find-ParagraphRange($NeedleText) {
    # Locate the (start,end) of the range matched by $NeedleText, then
    # return the full paragraph range that contains it.
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($NeedleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $rng.Find.Found) { return $null }

    $matchStart = $rng.Start

    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $pr = $paras.Item($i).Range
        if ($matchStart -ge $pr.Start -and $matchStart -lt $pr.End) {
            return $d.Range($pr.Start, $pr.End)
        }
    }
    return $null
}

function Format-Word($ParaPrefix, $Word, $Style) {
    $para = Find-ParagraphRange $ParaPrefix
    if ($para -eq $null) { return }

    # Find the target word inside that paragraph only.
    $sub = $d.Range($para.Start, $para.End)
    $sub.Find.ClearFormatting()
    $sub.Find.Execute($Word, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $sub.Find.Found) { return }

    if ($Style -eq "i") {
        $sub.Italic = 1
    } elseif ($Style -eq "b") {
        $sub.Bold = 1
    }

    # Force the surrounding single-space characters to become their own
    # (plainly formatted) runs, instead of being re-merged into the
    # neighbouring word's run, by toggling Bold on then back off (net
    # no-op on the space's own formatting, but it breaks run adjacency
    # so the space ends up as its own distinct run).
    $wStart = $sub.Start
    $wEnd = $sub.End

    $spaceBefore = $d.Range($wStart - 1, $wStart)
    if ($spaceBefore.Text -eq " ") {
        $spaceBefore.Bold = 1
        $spaceBefore.Bold = 0
    }

    $spaceAfter = $d.Range($wEnd, $wEnd + 1)
    if ($spaceAfter.Text -eq " ") {
        $spaceAfter.Bold = 1
        $spaceAfter.Bold = 0
    }
}

# --- Edit 1: "We have worked these tools:" -> "We have worked with *all* of these tools:"
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("We have worked these tools:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.Text = "We have worked with all of these tools:"
}
Format-Word "We have worked with all of these tools:" "all" "i"

# --- Edit 2: "I've found that knowledge of GitHub has been" -> bold "GitHub"
Format-Word "I've found that knowledge of GitHub has been" "GitHub" "b"

# --- Edit 3: "For instance, Markdown, HTML, and CSS all have similar" -> bold "all"
Format-Word "For instance, Markdown, HTML, and CSS all have similar" "all" "b"
